$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.916.17'
$ws.Range('D3').Value = '2.297.51'
$ws.Range('E3').Value = '  -1.02%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '299.94'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.44'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.515'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.55%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.507'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.75'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.87%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0787'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.78%  '
$ws.Range('E12').Value = '  +0.60%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '17.68'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.76'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.47%  '
$ws.Range('D15').Value = '2.656.56'
$ws.Range('E15').Value = '  -0.97%  '
$ws.Range('D16').Value = '2.264.03'
$ws.Range('E16').Value = '  -2.07%  '
$ws.Range('E17').Value = '  -2.41%  '
$ws.Range('D18').Value = '42.880.59'
$ws.Range('E18').Value = '  -0.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.65'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.16%  '
$ws.Range('D20').Value = '0.0₃0906'
$ws.Range('E20').Value = '  -0.81%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.09'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.00'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '241.05'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.13'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.32%  '
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('E27').Value = '  -0.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.91'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.52%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '166.56'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.08%  '
$ws.Range('E30').Value = '  -0.49%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.04'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.63%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '32.79'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.42%  '
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('B34').Value = 'RenderToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.75'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.41%  '
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.27%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.50'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.11%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.39'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.54%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0686'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.87%  '
$ws.Range('E39').Value = '  -2.21%  '
$ws.Range('E40').Value = '  -3.33%  '
$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.110'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.40%  '
$ws.Range('B42').Value = 'LidoDAOToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.75'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.00%  '
$ws.Range('D43').Value = '2.002.27'
$ws.Range('E43').Value = '  +0.29%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0284'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.74%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.12'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.29%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.14'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '17.26'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.22%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.77'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.78%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '2.522.63'
$ws.Range('E49').Value = '  -1.00%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '53.27'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.28%  '
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.15'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.99%  '
